$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Rename inline picture "name" attributes for the Pearson logo shown in the
# default (primary) footer and the first-page footer, and for the BTEC logo
# shown in the first-page header.
#
# InlineShape has no settable Name in the Word object model, so the
# documented approach is to convert the inline picture to a floating shape
# (which exposes a writable Name), rename it, then convert it back to an
# inline picture in place.

function Rename-InlinePicture($range, $newName) {
    $count = $range.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $inlineShp = $range.InlineShapes($i)
        $floatShp = $inlineShp.ConvertToShape()
        $floatShp.Name = $newName
        $floatShp.ConvertToInlineShape() | Out-Null
    }
}

# Default (primary) footer -> footer2.xml : PearsonLogo id=2, image1.png -> image2.png
Rename-InlinePicture $sec.Footers(1).Range "image2.png"

# First-page footer -> footer1.xml : PearsonLogo id=3, image1.png -> image2.png
Rename-InlinePicture $sec.Footers(2).Range "image2.png"

# First-page header -> header1.xml : BTec_Logo id=1, image2.jpg -> image1.jpg
Rename-InlinePicture $sec.Headers(2).Range "image1.jpg"
